$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 86: correct the date/time value in A86 ---
$ws.Range("A86").Value = 45461.2916666667

# --- Append new row 87 with the latest results from the R script ---

# Copy the date cell's style (number format "yyyy-mm-dd hh:mm:ss") from A86
# down to A87 before setting its value, so it renders/serializes the same way.
$ws.Range("A86").Copy($ws.Range("A87"))
$ws.Range("A87").Value = 45462.5994907407

$ws.Range("B87").Value = 3300
$ws.Range("C87").Value = 6.05999994277954
$ws.Range("D87").Value = 6
$ws.Range("E87").Value = 6.05999994277954
$ws.Range("F87").Value = 6

# G (adj_close) is stored as text in this sheet, same as the rest of column G.
$ws.Range("G87").NumberFormat = "@"
$ws.Range("G87").Value = "6"
$ws.Range("G87").Style = $ws.Range("G86").Style

# H (ticker) is a plain text value.
$ws.Range("H87").Value = "PAL.MI"
